$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("G").Delete()
